$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new data rows (rows 102-146) following the existing table pattern:
#   Column A (regcntr_id) cycles 10002..10010
#   Column B (device_id) increments sequentially starting at 3000121
#   Columns C-F are constant: "eng", TRUE, "superadmin", "now()"
$startRow = 102
$startDeviceId = 3000121
$count = 45

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $regCenterId = 10002 + ($i % 9)
    $deviceId = $startDeviceId + $i

    $ws.Cells.Item($row, 1).Value = $regCenterId
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Scroll the view so the newly added rows are visible and select the new range,
# matching the author's view state after pasting the rows in.
$ws.Application.ActiveWindow.ScrollRow = 128
$ws.Range("A102:F146").Select() | Out-Null

# Page setup: saved as portrait with an explicit print-quality (matches the
# pageSetup element emitted by Excel when the workbook is saved again).
$ws.PageSetup.Orientation = 1

